# This script re-orders the data rows of the "Artfynd" sheet.
# The underlying records (rows 2-13 and rows 20-24) were resorted; no values
# were actually added or removed, the rows were just rearranged.  We read the
# full row blocks into memory first (so the read is unaffected by later
# writes), then write every row back out in its new position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 2-13 (records reordered by their "Id" in column A) ---
$range1 = $ws.Range("A2:AY13")
$data1 = $range1.Value2

# target relative row (1-based, row 2 = 1) -> source relative row (1-based)
$map1 = @{
    1  = 10   # row 2  <- row 11
    2  = 11   # row 3  <- row 12
    3  = 6    # row 4  <- row 7
    4  = 7    # row 5  <- row 8
    5  = 1    # row 6  <- row 2
    6  = 8    # row 7  <- row 9
    7  = 2    # row 8  <- row 3
    8  = 12   # row 9  <- row 13
    9  = 3    # row 10 <- row 4
    10 = 4    # row 11 <- row 5
    11 = 5    # row 12 <- row 6
    12 = 9    # row 13 <- row 10
}

$cols1 = $data1.GetLength(1)
$newData1 = New-Object 'object[,]' 12, $cols1
for ($t = 1; $t -le 12; $t++) {
    $s = $map1[$t]
    for ($c = 1; $c -le $cols1; $c++) {
        $newData1[$t - 1, $c - 1] = $data1[$s, $c]
    }
}

# Columns Y and AA hold dates formatted as plain text (e.g. "2023-03-08").
# Force those columns to Text format before the bulk write so Excel does not
# "helpfully" reinterpret the strings as real date serials.
$ws.Range("Y2:Y13").NumberFormat = "@"
$ws.Range("AA2:AA13").NumberFormat = "@"

$range1.Value = $newData1

# Restore the default (unformatted) look of the cells now that the text
# value is safely stored, matching the rest of the sheet.
$ws.Range("Y2:Y13").ClearFormats()
$ws.Range("AA2:AA13").ClearFormats()

# --- Block 2: rows 20-24 (records reordered by their "Id" in column A) ---
$range2 = $ws.Range("A20:AY24")
$data2 = $range2.Value2

$map2 = @{
    1 = 3   # row 20 <- row 22
    2 = 1   # row 21 <- row 20
    3 = 4   # row 22 <- row 23
    4 = 5   # row 23 <- row 24
    5 = 2   # row 24 <- row 21
}

$cols2 = $data2.GetLength(1)
$newData2 = New-Object 'object[,]' 5, $cols2
for ($t = 1; $t -le 5; $t++) {
    $s = $map2[$t]
    for ($c = 1; $c -le $cols2; $c++) {
        $newData2[$t - 1, $c - 1] = $data2[$s, $c]
    }
}

$ws.Range("Y20:Y24").NumberFormat = "@"
$ws.Range("AA20:AA24").NumberFormat = "@"

$range2.Value = $newData2

$ws.Range("Y20:Y24").ClearFormats()
$ws.Range("AA20:AA24").ClearFormats()
